$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# -----------------------------------------------------------------
# Re-add "cluster 9" data: 5 new rows (78-82) at the bottom of the
# references sheet, each styled like a little mini "table block"
# (medium borders, Arial 10pt, a hyperlink-look column D) plus a
# newly restored hyperlink relationship on D73.
# -----------------------------------------------------------------

# Column A values (sail class)
$colA = @{
    78 = "Ilca 7"
    79 = "Ilca 6"
    80 = "49er"
    81 = "Ilca 7"
    82 = "Ilca 6"
}

# Column B values (competition name)
$colB = @{
    78 = "Semaine Olympique Francaise De Voile 2017"
    79 = "Semaine Olympique Francaise De Voile 2017"
    80 = "World Championship 2019"
    81 = "Semaine Olympique Francaise De Voile 2018"
    82 = "Semaine Olympique Francaise De Voile 2018"
}

# Column D values (source URL, stored as plain text - these are not
# live hyperlinks, just styled to look like the hyperlink column)
$colD = @{
    78 = "https://www.manage2sail.com/api/event/0adf7bcd-01d0-4214-a295-bb0b9136999e/regattaresult/056dd04a-3fad-45f1-b9c2-894bdb176b43"
    79 = "https://www.manage2sail.com/api/event/0adf7bcd-01d0-4214-a295-bb0b9136999e/regattaresult/056dd04a-3fad-45f1-b9c2-894bdb176b43"
    80 = "https://www.manage2sail.com/api/event/41de110b-ec26-427d-81fb-be7807677326/regattaresult/5c147b15-f550-4046-ab46-203d9390ddb7"
    81 = "https://www.manage2sail.com/api/event/71c3d3a9-60fc-4465-816d-4b474c3ef34b/regattaresult/f9177157-d702-4601-b6d2-1a12e644a0fb"
    82 = "https://www.manage2sail.com/api/event/71c3d3a9-60fc-4465-816d-4b474c3ef34b/regattaresult/f9177157-d702-4601-b6d2-1a12e644a0fb"
}

# Row heights (custom, as in the source workbook)
$rowHeights = @{
    78 = 18
    79 = 22.2
    80 = 21.6
    81 = 17.4
    82 = 23.4
}

foreach ($r in 78..82) {
    $ws.Rows($r).RowHeight = $rowHeights[$r]

    # --- Column A cell ---------------------------------------------------
    $a = $ws.Cells.Item($r, 1)
    $a.Value = $colA[$r]
    $a.WrapText = $true
    $a.Font.Name = "Arial"
    $a.Font.Size = 10
    $a.Font.Color = 0
    $a.Borders.LineStyle = 1
    $a.Borders.Weight = -4138
    $a.Borders.Color = 13421772
    $a.Borders.Item(7).Color = 0
    $a.Borders.Item(10).Color = 0

    # --- Column B cell -----------------------------------------------------
    $b = $ws.Cells.Item($r, 2)
    $b.Value = $colB[$r]
    $b.WrapText = $true
    $b.Font.Name = "Arial"
    $b.Font.Size = 10
    $b.Font.Color = 0
    $b.Borders.LineStyle = 1
    $b.Borders.Weight = -4138
    $b.Borders.Color = 13421772
    $b.Borders.Item(10).Color = 0
    if ($r -eq 80) {
        # this row's B cell also got a white fill
        $b.Interior.Color = 16777215
    }

    # --- Column D cell (hyperlink look-alike) -----------------------------
    $d = $ws.Cells.Item($r, 4)
    $d.Value = $colD[$r]
    $d.WrapText = $true
    $d.Font.Underline = 2
    $d.Font.ThemeColor = 11
    $d.Interior.Color = 16711935
    $d.Borders.LineStyle = 1
    $d.Borders.Weight = -4138
    $d.Borders.Color = 13421772
}

# -----------------------------------------------------------------
# Restore the hyperlink relationship on D73 (plain text value kept,
# just add the clickable link back on top of it).
# -----------------------------------------------------------------
$ws.Hyperlinks.Add($ws.Range("D73"), "https://manage2sail.com/api/event/d1e2dc90-a4bd-4065-9edb-36391f72670f/regattaresult/ae2cc32f-cc74-44dc-ba0a-51d395065fdc")

# -----------------------------------------------------------------
# A handful of rows above lost their manual 28.8pt row height -
# Excel auto-fit them back down when the sheet was re-saved.
# -----------------------------------------------------------------
foreach ($r in @(36, 44, 45, 46, 47, 50, 51, 55, 59, 60)) {
    $ws.Rows($r).AutoFit()
}

# Put the selection where the author left off
$ws.Range("D80").Select()
